# edit.ps1 -- apply "Updated cryptos list" refresh (Price / Volume(1h) columns,
# plus the 3-row reorder of Kaspa / dogwifhat / EthereumClassic and
# MantraDAO / USDe) to match the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = '90.419.61'
$ws.Cells.Item(2, 5).Value = '  -0.08%  '

# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = '3.090.67'
$ws.Cells.Item(3, 5).Value = '  -1.91%  '

# Row 4: TetherUSD
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  -0.12%  '

# Row 5: Solana
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '233.96'
$ws.Cells.Item(5, 5).Value = '  +8.98%  '

# Row 6: BNB
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '624.29'
$ws.Cells.Item(6, 5).Value = '  +0.12%  '

# Row 7: XRP
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.08'
$ws.Cells.Item(7, 5).Value = '  -4.16%  '

# Row 8: Dogecoin
$ws.Cells.Item(8, 5).Value = '  -0.54%  '

# Row 9: USDC
$ws.Cells.Item(9, 5).Value = '  +0.06%  '

# Row 10: LidoStakedEther
$ws.Cells.Item(10, 4).Value = '3.088.86'
$ws.Cells.Item(10, 5).Value = '  -1.85%  '

# Row 11: Cardano
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.727'
$ws.Cells.Item(11, 5).Value = '  -6.00%  '

# Row 12: TRON
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.196'
$ws.Cells.Item(12, 5).Value = '  -1.64%  '

# Row 13: Avalanche
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '36.65'
$ws.Cells.Item(13, 5).Value = '  +4.88%  '

# Row 14: ShibaInu
$ws.Cells.Item(14, 5).Value = '  +3.92%  '

# Row 15: Toncoin
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '5.46'
$ws.Cells.Item(15, 5).Value = '  -3.49%  '

# Row 16: WrappedBTC
$ws.Cells.Item(16, 4).Value = '90.049.40'
$ws.Cells.Item(16, 5).Value = '  -0.41%  '

# Row 17: WrappedliquidstakedEther2.0
$ws.Cells.Item(17, 5).Value = '  -2.15%  '

# Row 18: WrappedEther
$ws.Cells.Item(18, 4).Value = '3.065.72'
$ws.Cells.Item(18, 5).Value = '  -3.53%  '

# Row 19: SuiNetwork
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '3.77'
$ws.Cells.Item(19, 5).Value = '  +3.01%  '

# Row 20: PEPE
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.0000215'
$ws.Cells.Item(20, 5).Value = '  +2.15%  '

# Row 21: Chainlink
$ws.Cells.Item(21, 5).Value = '  -1.89%  '

# Row 22: BitcoinCash
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '436.25'
$ws.Cells.Item(22, 5).Value = '  -5.06%  '

# Row 23: Polkadot
$ws.Cells.Item(23, 5).Value = '  +6.66%  '

# Row 24: Uniswap
$ws.Cells.Item(24, 5).Value = '  -0.43%  '

# Row 25: NEARProtocol
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '5.93'
$ws.Cells.Item(25, 5).Value = '  +0.82%  '

# Row 26: LEO
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '7.58'
$ws.Cells.Item(26, 5).Value = '  -1.29%  '

# Row 27: Litecoin
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '88.49'
$ws.Cells.Item(27, 5).Value = '  -2.18%  '

# Row 28: Aptos
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '12.17'
$ws.Cells.Item(28, 5).Value = '  +0.90%  '

# Row 29: WrappedeETH
$ws.Cells.Item(29, 5).Value = '  -2.13%  '

# Row 30: Dai
$ws.Cells.Item(30, 5).Value = '  +0.04%  '

# Row 31: InternetComputer(DFINITY)
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '9.43'
$ws.Cells.Item(31, 5).Value = '  +2.58%  '

# Row 32: Cronos
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.158'
$ws.Cells.Item(32, 5).Value = '  -1.72%  '

# Row 33: Stellar
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.196'
$ws.Cells.Item(33, 5).Value = '  +9.25%  '

# Row 34: Kaspa
$ws.Cells.Item(34, 2).Value = 'Kaspa'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.154'
$ws.Cells.Item(34, 5).Value = '  +7.84%  '

# Row 35: dogwifhat
$ws.Cells.Item(35, 2).Value = 'dogwifhat'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '3.86'
$ws.Cells.Item(35, 5).Value = '  +6.03%  '

# Row 36: EthereumClassic
$ws.Cells.Item(36, 2).Value = 'EthereumClassic'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '25.77'
$ws.Cells.Item(36, 5).Value = '  -5.19%  '

# Row 37: Binance-PegBSC-USD
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.892'
$ws.Cells.Item(37, 5).Value = '  -11.06%  '

# Row 38: Bittensor
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '505.63'
$ws.Cells.Item(38, 5).Value = '  -1.71%  '

# Row 39: RenderToken
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '7.09'
$ws.Cells.Item(39, 5).Value = '  +3.33%  '

# Row 40: PancakeSwap
$ws.Cells.Item(40, 5).Value = '  +0.13%  '

# Row 41: Fetch.AI
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.29'
$ws.Cells.Item(41, 5).Value = '  +0.50%  '

# Row 42: Hedera
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.0869'
$ws.Cells.Item(42, 5).Value = '  -0.62%  '

# Row 43: WhiteBITCoin
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '22.19'
$ws.Cells.Item(43, 5).Value = '  -0.06%  '

# Row 44: PolygonEcosystemToken
$ws.Cells.Item(44, 5).Value = '  -1.92%  '

# Row 45: USDe
$ws.Cells.Item(45, 2).Value = 'USDe'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.00'
$ws.Cells.Item(45, 5).Value = '  -0.02%  '

# Row 46: MantraDAO
$ws.Cells.Item(46, 2).Value = 'MantraDAO'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '3.43'
$ws.Cells.Item(46, 5).Value = '  +53.05%  '

# Row 47: Stacks
$ws.Cells.Item(47, 5).Value = '  -1.07%  '

# Row 48: Monero
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '151.25'
$ws.Cells.Item(48, 5).Value = '  +2.83%  '

# Row 49: ARBITRUM
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.687'
$ws.Cells.Item(49, 5).Value = '  +4.18%  '

# Row 50: OKB
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '44.92'
$ws.Cells.Item(50, 5).Value = '  +1.08%  '

# Row 51: ImmutableX
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.34'
$ws.Cells.Item(51, 5).Value = '  +1.30%  '
